$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - Magnesium chloride unit price
$ws.Range("E8").Value2 = 0.38
$ws.Range("G8").Value2 = 0.349
$ws.Range("I8").Value2 = 0.411

# Row 9 - Zinc sulfate unit price
$ws.Range("E9").Value2 = 0.795
$ws.Range("G9").Value2 = 0.657
$ws.Range("I9").Value2 = 0.931

# Re-enter the Q8 formula so Excel regroups it into the main shared formula
# range (Q4:Q27) instead of keeping its own separate shared-formula group.
$ws.Range("Q8").Formula = "=IF(E8=H8, 1, IF(F8=`$F`$2, 1, 0))"

# Update selection to reflect the range the user selected while editing
$ws.Range("A8:XFD9").Select()
